$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 4222.1
$ws.Cells.Item(106, 9).Value = 3652.5
$ws.Cells.Item(106, 11).Value = 3652.5
$ws.Cells.Item(106, 13).Value = -3021.5

$ws.Cells.Item(107, 8).Value = 462.45
$ws.Cells.Item(107, 9).Value = 434.1579
$ws.Cells.Item(107, 10).Value = 1000
$ws.Cells.Item(107, 11).Value = 434.1579
$ws.Cells.Item(107, 12).Value = 1000
$ws.Cells.Item(107, 13).Value = 1485.8421
$ws.Cells.Item(107, 14).Value = -4840

$ws.Cells.Item(112, 8).Value = 3291.25
$ws.Cells.Item(112, 10).Value = 3828.4614
$ws.Cells.Item(112, 12).Value = 11485.3842
$ws.Cells.Item(112, 14).Value = -13701.3842

$ws.Cells.Item(137, 8).Value = 1843.0392
$ws.Cells.Item(137, 9).Value = 1868.8125
$ws.Cells.Item(137, 10).Value = 1799.6316
$ws.Cells.Item(137, 11).Value = 5606.4375
$ws.Cells.Item(137, 12).Value = 5398.8948
$ws.Cells.Item(137, 13).Value = -3056.4375
$ws.Cells.Item(137, 14).Value = -10498.8948

$ws.Cells.Item(138, 8).Value = 4039.9556
$ws.Cells.Item(138, 10).Value = 9482.3125
$ws.Cells.Item(138, 12).Value = 28446.9375
$ws.Cells.Item(138, 14).Value = -38726.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4749.3076
$ws.Cells.Item(32, 9).Value = 5426.878
$ws.Cells.Item(32, 11).Value = 5426.878
$ws.Cells.Item(32, 13).Value = -5139.878

$ws.Cells.Item(61, 8).Value = 2329.7144
$ws.Cells.Item(61, 9).Value = 2577.5
$ws.Cells.Item(61, 10).Value = 1999.3334
$ws.Cells.Item(61, 11).Value = 2577.5
$ws.Cells.Item(61, 12).Value = 1999.3334
$ws.Cells.Item(61, 13).Value = -2365.5
$ws.Cells.Item(61, 14).Value = -2423.3334

$ws.Cells.Item(74, 8).Value = 6257.926
$ws.Cells.Item(74, 9).Value = 8301.529
$ws.Cells.Item(74, 10).Value = 2783.8
$ws.Cells.Item(74, 11).Value = 8301.529
$ws.Cells.Item(74, 12).Value = 2783.8
$ws.Cells.Item(74, 13).Value = -7427.529
$ws.Cells.Item(74, 14).Value = -4531.8

$ws.Cells.Item(77, 8).Value = 6257.926
$ws.Cells.Item(77, 9).Value = 8301.529
$ws.Cells.Item(77, 10).Value = 2783.8
$ws.Cells.Item(77, 11).Value = 41507.645
$ws.Cells.Item(77, 12).Value = 13919
$ws.Cells.Item(77, 13).Value = -37139.645
$ws.Cells.Item(77, 14).Value = -22655

$ws.Cells.Item(132, 8).Value = 9526.536
$ws.Cells.Item(132, 9).Value = 8170.6875
$ws.Cells.Item(132, 10).Value = 11334.333
$ws.Cells.Item(132, 11).Value = 24512.0625
$ws.Cells.Item(132, 12).Value = 34002.999
$ws.Cells.Item(132, 13).Value = -21982.0625
$ws.Cells.Item(132, 14).Value = -39062.999

$ws.Cells.Item(136, 8).Value = 2329.7144
$ws.Cells.Item(136, 9).Value = 2577.5
$ws.Cells.Item(136, 10).Value = 1999.3334
$ws.Cells.Item(136, 11).Value = 7732.5
$ws.Cells.Item(136, 12).Value = 5998.0002
$ws.Cells.Item(136, 13).Value = -5182.5
$ws.Cells.Item(136, 14).Value = -11098.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 10110.75
$ws.Cells.Item(134, 9).Value = 10787.385
$ws.Cells.Item(134, 10).Value = 9647.789000000001
$ws.Cells.Item(134, 11).Value = 32362.155
$ws.Cells.Item(134, 12).Value = 28943.367
$ws.Cells.Item(134, 13).Value = -29827.155
$ws.Cells.Item(134, 14).Value = -34013.367

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7354621
$ws.Cells.Item(31, 9).Value = 1386.1951
$ws.Cells.Item(31, 10).Value = 18520644
$ws.Cells.Item(31, 11).Value = 1386.1951
$ws.Cells.Item(31, 12).Value = 18520644
$ws.Cells.Item(31, 13).Value = -1091.1951
$ws.Cells.Item(31, 14).Value = -18521234

$ws.Cells.Item(34, 8).Value = 7354621
$ws.Cells.Item(34, 9).Value = 1386.1951
$ws.Cells.Item(34, 10).Value = 18520644
$ws.Cells.Item(34, 11).Value = 1386.1951
$ws.Cells.Item(34, 12).Value = 18520644
$ws.Cells.Item(34, 13).Value = -1184.1951
$ws.Cells.Item(34, 14).Value = -18521048

$ws.Cells.Item(132, 8).Value = 3228.0322
$ws.Cells.Item(132, 9).Value = 2744.3845
$ws.Cells.Item(132, 10).Value = 3577.3333
$ws.Cells.Item(132, 11).Value = 8233.1535
$ws.Cells.Item(132, 12).Value = 10731.9999
$ws.Cells.Item(132, 13).Value = -5703.1535
$ws.Cells.Item(132, 14).Value = -15791.9999

$ws.Cells.Item(134, 8).Value = 3490.1052
$ws.Cells.Item(134, 9).Value = 1827.1428
$ws.Cells.Item(134, 10).Value = 4460.1665
$ws.Cells.Item(134, 11).Value = 5481.428400000001
$ws.Cells.Item(134, 12).Value = 13380.4995
$ws.Cells.Item(134, 13).Value = -2946.428400000001
$ws.Cells.Item(134, 14).Value = -18450.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(104, 8).Value = 8336083.5
$ws.Cells.Item(104, 9).Value = 100000000
$ws.Cells.Item(104, 10).Value = 3000
$ws.Cells.Item(104, 11).Value = 300000000
$ws.Cells.Item(104, 12).Value = 9000
$ws.Cells.Item(104, 13).Value = -299997379
$ws.Cells.Item(104, 14).Value = -14242

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7285.421
$ws.Cells.Item(132, 9).Value = 13803.125
$ws.Cells.Item(132, 10).Value = 2545.2727
$ws.Cells.Item(132, 11).Value = 41409.375
$ws.Cells.Item(132, 12).Value = 7635.8181
$ws.Cells.Item(132, 13).Value = -38879.375
$ws.Cells.Item(132, 14).Value = -12695.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3822.5
$ws.Cells.Item(7, 9).Value = 3433.3333
$ws.Cells.Item(7, 10).Value = 4056
$ws.Cells.Item(7, 11).Value = 3433.3333
$ws.Cells.Item(7, 12).Value = 4056
$ws.Cells.Item(7, 13).Value = -3321.3333
$ws.Cells.Item(7, 14).Value = -4280

$ws.Cells.Item(40, 8).Value = 3844.4666
$ws.Cells.Item(40, 9).Value = 3304.75
$ws.Cells.Item(40, 10).Value = 6003.3335
$ws.Cells.Item(40, 11).Value = 3304.75
$ws.Cells.Item(40, 12).Value = 6003.3335
$ws.Cells.Item(40, 13).Value = -3168.75
$ws.Cells.Item(40, 14).Value = -6275.3335

$ws.Cells.Item(126, 8).Value = 3822.5
$ws.Cells.Item(126, 9).Value = 3433.3333
$ws.Cells.Item(126, 10).Value = 4056
$ws.Cells.Item(126, 11).Value = 10299.9999
$ws.Cells.Item(126, 12).Value = 12168
$ws.Cells.Item(126, 13).Value = -7829.999899999999
$ws.Cells.Item(126, 14).Value = -17108

$ws.Cells.Item(132, 8).Value = 43482230
$ws.Cells.Item(132, 9).Value = 58826020
$ws.Cells.Item(132, 10).Value = 8167.5
$ws.Cells.Item(132, 11).Value = 176478060
$ws.Cells.Item(132, 12).Value = 24502.5
$ws.Cells.Item(132, 13).Value = -176475530
$ws.Cells.Item(132, 14).Value = -29562.5

$ws.Cells.Item(136, 8).Value = 41668056
$ws.Cells.Item(136, 9).Value = 55556856
$ws.Cells.Item(136, 10).Value = 1652.6666
$ws.Cells.Item(136, 11).Value = 166670568
$ws.Cells.Item(136, 12).Value = 4957.9998
$ws.Cells.Item(136, 13).Value = -166668018
$ws.Cells.Item(136, 14).Value = -10057.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1145
$ws.Cells.Item(81, 9).Value = 800
$ws.Cells.Item(81, 10).Value = 1260
$ws.Cells.Item(81, 11).Value = 1600
$ws.Cells.Item(81, 12).Value = 2520
$ws.Cells.Item(81, 13).Value = -539
$ws.Cells.Item(81, 14).Value = -4642

$ws.Cells.Item(84, 8).Value = 1145
$ws.Cells.Item(84, 9).Value = 800
$ws.Cells.Item(84, 10).Value = 1260
$ws.Cells.Item(84, 11).Value = 8000
$ws.Cells.Item(84, 12).Value = 12600
$ws.Cells.Item(84, 13).Value = -2696
$ws.Cells.Item(84, 14).Value = -23208

$ws.Cells.Item(122, 8).Value = 3775.6667
$ws.Cells.Item(122, 9).Value = 4661.3447
$ws.Cells.Item(122, 10).Value = 1207.2
$ws.Cells.Item(122, 11).Value = 13984.0341
$ws.Cells.Item(122, 12).Value = 3621.6
$ws.Cells.Item(122, 13).Value = -11534.0341
$ws.Cells.Item(122, 14).Value = -8521.6

$ws.Cells.Item(132, 8).Value = 3663.44
$ws.Cells.Item(132, 9).Value = 4441.143
$ws.Cells.Item(132, 10).Value = 3361
$ws.Cells.Item(132, 11).Value = 13323.429
$ws.Cells.Item(132, 12).Value = 10083
$ws.Cells.Item(132, 13).Value = -10793.429
$ws.Cells.Item(132, 14).Value = -15143
